$d = $word.ActiveDocument

# --- Helper: simple literal Find & Replace over the whole document content ---
function Replace-Text($find, $replace) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Host "WARNING: not found -> $find"
    }
    return $ok
}

# 1) "Noticias actualizadas sobre la liga" -> add trailing "."
Replace-Text "Noticias actualizadas sobre la liga" "Noticias actualizadas sobre la liga."

# 2) "...podrán recibir recomendación sobre fichajes a realizar" -> "...podrán recibir recomendaciones sobre fichajes a realizar"
Replace-Text "podrán recibir recomendación sobre fichajes a realizar" "podrán recibir recomendaciones sobre fichajes a realizar"

# 3) "Jugadores similares a otros según su rendimiento actual" -> add trailing "."
Replace-Text "Jugadores similares a otros según su rendimiento actual" "Jugadores similares a otros según su rendimiento actual."

# 4) "(Visitante. Administrado y Usuario registrado)" -> "(Visitante. Administrador y Usuario registrado)"
Replace-Text "Visitante. Administrado y Usuario registrado" "Visitante. Administrador y Usuario registrado"

# 5) "...de los usuarios de las Liga Fantasy." -> "...de los usuarios de las Ligas Fantasy."
Replace-Text "de los usuarios de las Liga Fantasy." "de los usuarios de las Ligas Fantasy."

# 6) "...como son tratado de las urls, settings, admin, etc." -> "...como son el tratado de las urls, settings, admin, etc."
Replace-Text "como son tratado de las urls, settings, admin, etc." "como son el tratado de las urls, settings, admin, etc."

# 7) Bootstrap paragraph: "diseño y adopción de la aplicación web más profesional" -> "diseño y adaptación más profesional a la aplicación web"
Replace-Text "diseño y adopción de la aplicación web más profesional" "diseño y adaptación más profesional a la aplicación web"

# 8) "Junto con HTML5 y CSS3 componen" -> "Junto con HTML5, CSS3 y JavaScript, componen"
Replace-Text "Junto con HTML5 y CSS3 componen" "Junto con HTML5, CSS3 y JavaScript, componen"

Write-Host "Done"
